# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> linked from the (only) Slide Master  ("Integral")
#   ppt/theme/theme2.xml  -> linked from the Notes Master          ("Office Theme")
#
# The authored change swaps the two themes' content wholesale: theme1.xml
# ends up holding the "Office Theme" palette and theme2.xml ends up holding
# the "Integral" palette (font scheme / format scheme are identical between
# the two themes already, so only the 12-slot colour scheme - and the theme
# name - actually differ).
#
# The PowerPoint object model only exposes a single editable Theme object in
# this deck (reached via SlideMaster.Theme, NotesMaster.Theme, Design.*,
# HandoutMaster.Theme, Slide.ThemeColorScheme, ... they all resolve to the
# same underlying theme part backing the Slide Master / theme1.xml). There is
# no COM surface onto the Notes Master's own theme part, so we apply the
# reachable half of the swap: push the "Office Theme" colour scheme into the
# Slide Master's theme.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# ThemeColorScheme.Item index order, confirmed against the existing
# "Integral" values: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink.
$colorScheme.Item(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1      000000
$colorScheme.Item(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1      FFFFFF
$colorScheme.Item(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2      44546A
$colorScheme.Item(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2      E7E6E6
$colorScheme.Item(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1  5B9BD5
$colorScheme.Item(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2  ED7D31
$colorScheme.Item(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3  A5A5A5
$colorScheme.Item(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4  FFC000
$colorScheme.Item(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5  4472C4
$colorScheme.Item(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6  70AD47
$colorScheme.Item(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink    0563C1
$colorScheme.Item(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink 954F72
